# Q3 Update - 2025
# Applies the UN-SEN.xlsx Q3 data refresh:
#   1. Refreshes the shared "short-url" value used by every data row.
#   2. Updates a handful of refugees/asylum_seekers figures for existing rows.
#   3. Appends a new data row (847) for Yemen as a country of origin.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: assign a value while preserving text storage for numeric-looking
# strings (this workbook stores every cell - including numbers - as text).
function Set-TextSafeValue {
    param($range, [string]$value)
    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $range.NumberFormat = "@"
    }
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# 1. Update the short-url shared across every row (B2:B846 all reference the
#    same value, so re-assigning the whole range keeps that single shared
#    string in sync everywhere).
# ---------------------------------------------------------------------------
$ws.Range("B2:B846").Value = "nzPm4F"

# ---------------------------------------------------------------------------
# 2. Refresh refugees (N) / asylum_seekers (O) figures for the affected rows.
# ---------------------------------------------------------------------------
$rowUpdates = @{
    823 = @{ N = "252"; O = "509" }
    824 = @{ O = "22" }
    825 = @{ N = "0";   O = "17" }
    826 = @{ O = "40" }
    827 = @{ O = "35" }
    830 = @{ O = "7" }
    831 = @{ N = "17";  O = "125" }
    833 = @{ N = "7";   O = "59" }
    834 = @{ N = "6";   O = "10" }
    836 = @{ N = "11539"; O = "35" }
    839 = @{ O = "10" }
    840 = @{ O = "225" }
    841 = @{ O = "0" }
    842 = @{ O = "10" }
    843 = @{ O = "27" }
    845 = @{ O = "13" }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        Set-TextSafeValue $ws.Range($addr) $cols[$col]
    }
}

# ---------------------------------------------------------------------------
# 3. Append the new Yemen row (847). Copy the formatting from the row above
#    so styles/number formats line up with the rest of the table, then fill
#    in the values.
# ---------------------------------------------------------------------------
$ws.Range("A846:V846").Copy() | Out-Null
$ws.Range("A847:V847").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newRow = [ordered]@{
    A = "1"
    B = "nzPm4F"
    C = "1"
    D = "846"
    E = "2024"
    F = "211"
    G = "Yemen"
    H = "YEM"
    I = "YEM"
    J = "165"
    K = "Senegal"
    L = "SEN"
    M = "SEN"
    N = "0"
    O = "5"
    P = "0"
    Q = "0"
    R = "0"
    S = "0"
    T = "0"
    U = "-"
    V = "0"
}

foreach ($col in $newRow.Keys) {
    $addr = $col + "847"
    Set-TextSafeValue $ws.Range($addr) $newRow[$col]
}

Write-Output "Applied Q3 2025 update: short-url refresh, 16 row figure updates, and new Yemen row (847)."
